# Updated cryptos list - refresh Price (D) and Volume(1h) (E) columns.
# Values in column D that "look like" a plain number (single decimal point)
# must stay TEXT (as they were originally stored as inline strings), so for
# those we briefly force a Text number format, assign the value, then clear
# the formatting again so no stray style index is left behind on the cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($addr, $value)
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.ClearFormats()
}

# Row 2 - Bitcoin
Set-TextValue "D2" "68.217.99"
$ws.Range("E2").Value = "  +0.23%  "

# Row 3 - Ethereum
Set-TextValue "D3" "3.673.99"
$ws.Range("E3").Value = "  -3.37%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.09%  "

# Row 5 - BNB
Set-TextValue "D5" "596.66"
$ws.Range("E5").Value = "  +0.32%  "

# Row 6 - Solana
Set-TextValue "D6" "165.96"
$ws.Range("E6").Value = "  -3.97%  "

# Row 7 - LidoStakedEther
Set-TextValue "D7" "3.667.33"
$ws.Range("E7").Value = "  -3.59%  "

# Row 8 - USDC
$ws.Range("E8").Value = "  +0.05%  "

# Row 9 - XRP
$ws.Range("E9").Value = "  +0.88%  "

# Row 10 - Dogecoin
$ws.Range("E10").Value = "  +2.85%  "

# Row 11 - Toncoin
Set-TextValue "D11" "6.27"
$ws.Range("E11").Value = "  +0.15%  "

# Row 12 - Cardano
$ws.Range("E12").Value = "  -1.93%  "

# Row 13 - Avalanche
Set-TextValue "D13" "37.84"
$ws.Range("E13").Value = "  -0.53%  "

# Row 14 - ShibaInu
$ws.Range("E14").Value = "  -0.41%  "

# Row 15 - WrappedliquidstakedEther2.0
Set-TextValue "D15" "4.293.63"
$ws.Range("E15").Value = "  -3.20%  "

# Row 16 - WrappedEther
Set-TextValue "D16" "3.677.89"
$ws.Range("E16").Value = "  -3.10%  "

# Row 17 - WrappedBTC
Set-TextValue "D17" "68.084.27"
$ws.Range("E17").Value = "  -0.04%  "

# Row 18 - Polkadot
$ws.Range("E18").Value = "  +0.62%  "

# Row 19 - TRON
$ws.Range("E19").Value = "  -0.92%  "

# Row 20 - Chainlink
Set-TextValue "D20" "17.06"
$ws.Range("E20").Value = "  +6.11%  "

# Row 21 - BitcoinCash
Set-TextValue "D21" "490.81"
$ws.Range("E21").Value = "  +0.17%  "

# Row 22 - Uniswap
Set-TextValue "D22" "9.08"
$ws.Range("E22").Value = "  -2.23%  "

# Row 23 - Polygon
$ws.Range("E23").Value = "  -1.99%  "

# Row 24 - Litecoin
Set-TextValue "D24" "84.32"
$ws.Range("E24").Value = "  -0.38%  "

# Row 25 - PEPE
Set-TextValue "D25" "0.0000142"
$ws.Range("E25").Value = "  +3.00%  "

# Row 26 - Fetch.AI
$ws.Range("E26").Value = "  -4.57%  "

# Row 27 - InternetComputer(DFINITY)
Set-TextValue "D27" "12.16"

# Row 28 - RenderToken
Set-TextValue "D28" "10.02"
$ws.Range("E28").Value = "  -2.24%  "

# Row 29 - Dai
$ws.Range("E29").Value = "  -0.02%  "

# Row 30 - PancakeSwap
$ws.Range("E30").Value = "  -0.97%  "

# Row 31 - ImmutableX
Set-TextValue "D31" "2.37"
$ws.Range("E31").Value = "  -2.62%  "

# Row 32 - NEARProtocol
$ws.Range("E32").Value = "  +1.49%  "

# Row 33 - EthereumClassic
Set-TextValue "D33" "31.24"
$ws.Range("E33").Value = "  -4.79%  "

# Row 34 - WrappedeETH
Set-TextValue "D34" "3.816.62"
$ws.Range("E34").Value = "  -3.18%  "

# Row 35 - Hedera
$ws.Range("E35").Value = "  -1.50%  "

# Row 36 - RenzoRestakedETH
Set-TextValue "D36" "3.619.86"
$ws.Range("E36").Value = "  -3.22%  "

# Row 37 - FirstDigitalUSD
Set-TextValue "D37" "0.999"
$ws.Range("E37").Value = "  +0.03%  "

# Row 38 - Mantle
Set-TextValue "D38" "0.992"
$ws.Range("E38").Value = "  -1.82%  "

# Row 39 - Filecoin
Set-TextValue "D39" "5.72"
$ws.Range("E39").Value = "  -1.10%  "

# Row 40 - Kaspa
$ws.Range("E40").Value = "  -3.84%  "

# Row 41 - TheGraph
$ws.Range("E41").Value = "  -2.13%  "

# Row 42 - Bittensor
Set-TextValue "D42" "430.52"
$ws.Range("E42").Value = "  -5.49%  "

# Row 43 - OKB
Set-TextValue "D43" "48.64"
$ws.Range("E43").Value = "  -0.62%  "

# Row 44 - Stacks
Set-TextValue "D44" "1.95"
$ws.Range("E44").Value = "  -2.66%  "

# Row 45 - dogwifhat
$ws.Range("E45").Value = "  -3.03%  "

# Row 46 - Cosmos
Set-TextValue "D46" "8.35"
$ws.Range("E46").Value = "  +0.82%  "

# Row 47 - USDe
$ws.Range("E47").Value = "  +0.02%  "

# Row 48 - Arweave
Set-TextValue "D48" "40.24"
$ws.Range("E48").Value = "  -2.99%  "

# Row 49 - Monero
Set-TextValue "D49" "141.26"
$ws.Range("E49").Value = "  +1.31%  "

# Row 50 - Maker
Set-TextValue "D50" "2.726.49"
$ws.Range("E50").Value = "  -3.70%  "

# Row 51 - VeChain
$ws.Range("E51").Value = "  -1.09%  "
